$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.464.24"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.568.14"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'208.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'0.502"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "'22.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "'0.250"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.789.80"
$ws.Range("D13").Value = "1.599.43"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "'3.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "'63.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "27.464.60"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "'214.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").Value = "0.0₃0691"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "'7.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'4.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("D25").Value = "'152.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("D28").Value = "'15.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").Value = "'0.105"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("D32").Value = "'3.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "1.376.81"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D36").Value = "'0.959"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "'0.0168"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").Value = "'0.553"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("D40").Value = "'0.826"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Value = "'0.977"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'1.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.88%  "
$ws.Range("D44").Value = "'64.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("D45").Value = "'2.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'5.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "1.703.13"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").Value = "'0.0497"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
